# ---------------------------------------------------------------------------
# updating functions and inputs file to accept inputs from spreadsheet
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. settings: nr_patches value changes from 6 to 2
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B2").Value = 2

# ---------------------------------------------------------------------------
# 2. jurisdiction: drop PA / Germany / France / Italy, keep only NY & NJ
# ---------------------------------------------------------------------------
$jurisdiction = $wb.Worksheets.Item("jurisdiction")
$jurisdiction.Range("A4:D7").ClearContents()

# ---------------------------------------------------------------------------
# 3. parameters: completely new layout
#    parameter | baseline | description | min | max | distribution | source
# ---------------------------------------------------------------------------
$parameters = $wb.Worksheets.Item("parameters")
$parameters.Cells.ClearContents()

$parameters.Range("A1").Value = "parameter"
$parameters.Range("B1").Value = "baseline"
$parameters.Range("C1").Value = "description"
$parameters.Range("D1").Value = "min"
$parameters.Range("E1").Value = "max"
$parameters.Range("F1").Value = "distribution"
$parameters.Range("G1").Value = "source"

$parameters.Range("A2").Value = "sigma"
$parameters.Range("B2").Formula = "=1/6"

$parameters.Range("A3").Value = "delta"
$parameters.Range("B3").Formula = "=1/6"

$parameters.Range("A4").Value = "gamma"
$parameters.Range("B4").Formula = "=1/3"

$parameters.Range("A5").Value = "tau"
$parameters.Range("B5").Value = 0.15

$parameters.Range("A6").Value = "c"
$parameters.Range("B6").Value = 1

$parameters.Range("A7").Value = "obs_lag"
$parameters.Range("B7").Value = 5

$parameters.Range("A8").Value = "days_to_adjust_NPI"
$parameters.Range("B8").Value = 7

$parameters.Columns.Item(1).ColumnWidth = 16.666666666666668
$parameters.Columns.Item(3).ColumnWidth = 26.5

# ---------------------------------------------------------------------------
# 4. beta: brand-new sheet, identity-ish matrix keyed off jurisdiction
# ---------------------------------------------------------------------------
$beta = $wb.Worksheets.Add()
$beta.Name = "beta"

$beta.Range("A1").Value = "jurisdiction"
$beta.Range("B1").Formula = "=TRANSPOSE(A2:A3)"

$beta.Range("A2").Formula = "=jurisdiction!B2"
$beta.Range("B2").Value = 1
$beta.Range("C2").Value = 0

$beta.Range("A3").Formula = "=jurisdiction!B3"
$beta.Range("B3").Value = 0
$beta.Range("C3").Value = 1

# ---------------------------------------------------------------------------
# 5. reorder tabs: settings, parameters, jurisdiction, beta, travel, relative-mixing
#    (re-fetch sheet handles by name after each structural change - stale
#    references can rebind to the wrong sheet once Add()/Move() reshuffle
#    the worksheets collection)
# ---------------------------------------------------------------------------
$jurisdiction = $wb.Worksheets.Item("jurisdiction")
$parameters = $wb.Worksheets.Item("parameters")
$parameters.Move($jurisdiction)

$travel = $wb.Worksheets.Item("travel")
$beta = $wb.Worksheets.Item("beta")
$beta.Move($travel)

# ---------------------------------------------------------------------------
# 6. sheet-view bookkeeping (zoom / selection / active tab)
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
$excel.ActiveWindow.Zoom = 181
$settings.Range("B7").Select()

$jurisdiction = $wb.Worksheets.Item("jurisdiction")
$jurisdiction.Activate()
$jurisdiction.Range("A1:D8").Select()

$beta = $wb.Worksheets.Item("beta")
$beta.Activate()
$beta.Range("A4").Select()

$travel = $wb.Worksheets.Item("travel")
$travel.Activate()
$travel.Range("D1:G1").Select()

$relativeMixing = $wb.Worksheets.Item("relative-mixing")
$relativeMixing.Activate()
$relativeMixing.Range("G25").Select()

$parameters = $wb.Worksheets.Item("parameters")
$parameters.Activate()
$excel.ActiveWindow.Zoom = 141
$parameters.Range("C4").Select()
